# Auto-generated Excel COM-interop script replicating the Marilith_Profits diff.
# Each block corresponds to one row in one sheet that had numeric data refreshed
# (market-price / profit recompute). Cells that the diff adds get a new .Value;
# cells the diff drops entirely are cleared with ClearContents() so the cell node
# disappears from the OOXML, matching upstream.

$wb = $excel.ActiveWorkbook

# Sheet ALC, row 20 (hunk 0)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 990.1667
$ws.Range("I20").Value = 990.1667
$ws.Range("K20").Value = 990.1667
$ws.Range("M20").Value = -760.1667

# Sheet ALC, row 32 (hunk 1)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 714.2857
$ws.Range("I32").Value = 550
$ws.Range("K32").Value = 550
$ws.Range("M32").Value = -224

# Sheet ALC, row 35 (hunk 2)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H35").Value = 990.1667
$ws.Range("I35").Value = 990.1667
$ws.Range("K35").Value = 990.1667
$ws.Range("M35").Value = -611.1667

# Sheet ALC, row 51 (hunk 3)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4997
$ws.Range("J51").Value = 4997
$ws.Range("L51").Value = 4997
$ws.Range("N51").Value = -5965

# Sheet ALC, row 88 (hunk 4)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1184.5
$ws.Range("J88").Value = 1330.6364
$ws.Range("L88").Value = 1330.6364
$ws.Range("N88").Value = -2142.6364

# Sheet ALC, row 91 (hunk 5)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 1184.5
$ws.Range("J91").Value = 1330.6364
$ws.Range("L91").Value = 1330.6364
$ws.Range("N91").Value = -4138.6364

# Sheet ALC, row 100 (hunk 6)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J100").Value = 4000
$ws.Range("L100").Value = 4000
$ws.Range("N100").Value = -5082

# Sheet ALC, row 111 (hunk 7)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1127.5714
$ws.Range("I111").Value = 1127.5714
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 3382.7142
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -315.7142000000003
$ws.Range("N111").ClearContents()

# Sheet ALC, row 132 (hunk 8)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4330.8
$ws.Range("I132").Value = 4423.125
$ws.Range("K132").Value = 13269.375
$ws.Range("M132").Value = -10739.375

# Sheet ARM, row 61 (hunk 9)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3043.889
$ws.Range("I61").Value = 1723.75
$ws.Range("K61").Value = 1723.75
$ws.Range("M61").Value = -1511.75

# Sheet ARM, row 74 (hunk 10)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4550
$ws.Range("I74").Value = 4550
$ws.Range("K74").Value = 4550
$ws.Range("M74").Value = -3676

# Sheet ARM, row 77 (hunk 11)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4550
$ws.Range("I77").Value = 4550
$ws.Range("K77").Value = 22750
$ws.Range("M77").Value = -18382

# Sheet ARM, row 110 (hunk 12)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2258.6667
$ws.Range("I110").Value = 571.7692
$ws.Range("K110").Value = 571.7692
$ws.Range("M110").Value = 1473.2308

# Sheet ARM, row 136 (hunk 13)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3043.889
$ws.Range("I136").Value = 1723.75
$ws.Range("K136").Value = 5171.25
$ws.Range("M136").Value = -2621.25

# Sheet BSM, row 22 (hunk 14)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 256.2
$ws.Range("I22").Value = 326.66666
$ws.Range("J22").Value = 150.5
$ws.Range("K22").Value = 326.66666
$ws.Range("L22").Value = 150.5
$ws.Range("M22").Value = -153.66666
$ws.Range("N22").Value = -496.5

# Sheet BSM, row 106 (hunk 15)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 20223.666
$ws.Range("J106").Value = 20223.666
$ws.Range("L106").Value = 20223.666
$ws.Range("N106").Value = -22747.666

# Sheet BSM, row 134 (hunk 16)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4069.5833
$ws.Range("I134").Value = 4311
$ws.Range("J134").Value = 1414
$ws.Range("K134").Value = 12933
$ws.Range("L134").Value = 4242
$ws.Range("M134").Value = -10398
$ws.Range("N134").Value = -9312

# Sheet CRP, row 14 (hunk 17)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

# Sheet CRP, row 22 (hunk 18)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 900
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 900
$ws.Range("N22").Value = -1600
$ws.Range("M22").ClearContents()

# Sheet CRP, row 58 (hunk 19)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3577.9167
$ws.Range("I58").Value = 3195.6
$ws.Range("J58").Value = 3851
$ws.Range("K58").Value = 3195.6
$ws.Range("L58").Value = 3851
$ws.Range("M58").Value = -2992.6
$ws.Range("N58").Value = -4257

# Sheet CRP, row 94 (hunk 20)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2723.5
$ws.Range("I94").Value = 2648.6667
$ws.Range("K94").Value = 2648.6667
$ws.Range("M94").Value = -2197.6667

# Sheet CRP, row 99 (hunk 21)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1966.6666
$ws.Range("I99").Value = 1450
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 1450
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = 48
$ws.Range("N99").Value = -5996

# Sheet CRP, row 126 (hunk 22)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1966.6666
$ws.Range("I126").Value = 1450
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 4350
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -1880
$ws.Range("N126").Value = -13940

# Sheet CRP, row 134 (hunk 23)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3080.8
$ws.Range("I134").Value = 1562.7
$ws.Range("K134").Value = 4688.1
$ws.Range("M134").Value = -2153.1

# Sheet CRP, row 136 (hunk 24)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3577.9167
$ws.Range("I136").Value = 3195.6
$ws.Range("J136").Value = 3851
$ws.Range("K136").Value = 9586.799999999999
$ws.Range("L136").Value = 11553
$ws.Range("M136").Value = -7036.799999999999
$ws.Range("N136").Value = -16653

# Sheet CUL, row 98 (hunk 25)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 381.25
$ws.Range("J98").Value = 374.5
$ws.Range("L98").Value = 1123.5
$ws.Range("N98").Value = -4119.5

# Sheet GSM, row 136 (hunk 26)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 38598.8
$ws.Range("J136").Value = 38598.8
$ws.Range("L136").Value = 115796.4
$ws.Range("N136").Value = -120896.4

# Sheet LTW, row 22 (hunk 27)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 870.1
$ws.Range("J22").Value = 750.5
$ws.Range("L22").Value = 750.5
$ws.Range("N22").Value = -1340.5

# Sheet LTW, row 27 (hunk 28)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 870.1
$ws.Range("J27").Value = 750.5
$ws.Range("L27").Value = 750.5
$ws.Range("N27").Value = -964.5

# Sheet LTW, row 30 (hunk 29)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 415.83334
$ws.Range("I30").Value = 365
$ws.Range("K30").Value = 365
$ws.Range("M30").Value = -257

# Sheet LTW, row 40 (hunk 30)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6724.75
$ws.Range("I40").Value = 6589.7
$ws.Range("J40").Value = 7400
$ws.Range("K40").Value = 6589.7
$ws.Range("L40").Value = 7400
$ws.Range("M40").Value = -6453.7
$ws.Range("N40").Value = -7672

# Sheet LTW, row 54 (hunk 31)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H54").Value = 46563
$ws.Range("J54").Value = 46563
$ws.Range("L54").Value = 46563
$ws.Range("N54").Value = -47851

# Sheet LTW, row 68 (hunk 32)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5849.5
$ws.Range("I68").Value = 5966
$ws.Range("J68").Value = 5500
$ws.Range("K68").Value = 5966
$ws.Range("L68").Value = 5500
$ws.Range("M68").Value = -5217
$ws.Range("N68").Value = -6998

# Sheet LTW, row 71 (hunk 33)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 5849.5
$ws.Range("I71").Value = 5966
$ws.Range("J71").Value = 5500
$ws.Range("K71").Value = 29830
$ws.Range("L71").Value = 27500
$ws.Range("M71").Value = -26086
$ws.Range("N71").Value = -34988

# Sheet LTW, row 132 (hunk 34)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 11399.25
$ws.Range("I132").Value = 11146.111
$ws.Range("K132").Value = 33438.333
$ws.Range("M132").Value = -30908.333

# Sheet WVR, row 104 (hunk 35)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 9996
$ws.Range("J104").Value = 9996
$ws.Range("L104").Value = 9996
$ws.Range("N104").Value = -16984

# Sheet WVR, row 123 (hunk 36)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 20000
$ws.Range("J123").Value = 20000
$ws.Range("L123").Value = 20000
$ws.Range("N123").Value = -29800

# Sheet WVR, row 132 (hunk 37)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1328.2858
$ws.Range("I132").Value = 1328.2858
$ws.Range("K132").Value = 3984.8574
$ws.Range("M132").Value = -1454.8574
